$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. B16/C16: was "material_name"/"varchar" (plain box style). The row now
#    becomes a foreign-key-style row referencing a new "study_program" table,
#    so it gets the same box+fill style used by other id-reference rows
#    (e.g. B8 "rule_id"/"int").
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "study_program_id"
$ws.Range("C16").Value = "int"

$srcIdStyle = $ws.Range("B8:C8")
$srcIdStyle.Copy()
$ws.Range("B16:C16").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. B17/C17 (new row): gets the old "material_name"/"varchar" content that
#    used to live in B16/C16, in the regular plain box style (same as B2).
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "material_name"
$ws.Range("C17").Value = "varchar"

$srcPlainStyle = $ws.Range("B2:C2")
$srcPlainStyle.Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. B19:C19 (new merged header row): new mini-table header "study_program",
#    styled like the other table headers (fill + left align) but with a full
#    box border around each cell instead of the split left/right border.
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "study_program"

$srcHeaderStyle = $ws.Range("B1:C1")
$srcHeaderStyle.Copy()
$dstHeader = $ws.Range("B19:C19")
$dstHeader.PasteSpecial(-4122)

$ws.Range("B19:C19").Merge()

$dstHeader.Borders.LineStyle = -4142
$dstHeader.Borders.LineStyle = 1
$dstHeader.Borders.Weight = 2

# ---------------------------------------------------------------------------
# 4. B20/C20 (new row): regular field row "user_login"/"int" (plain box).
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "user_login"
$ws.Range("C20").Value = "int"

$srcPlainStyle.Copy()
$ws.Range("B20:C20").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. B21/C21 (new row): "study_program"/"varchar" field row, plain box style
#    (visually identical to the regular field rows).
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "study_program"
$ws.Range("C21").Value = "varchar"

$srcPlainStyle.Copy()
$ws.Range("B21:C21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. B22/C22 (new row): blank spacer cells with no border/fill.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""

# ---------------------------------------------------------------------------
# 7. Update selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("B10:C10").Select()
